$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.185402274131775
$ws.Range("B1").Value = 2.472376346588135
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.762502789497375
$ws.Range("E1").Value = 1.179947137832642
